$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoopFilter LPF")

# Update the scale factor value (B16): 2 -> 4
$ws.Range("B16").Value = 4

# Update the active selection to B17
$ws.Range("B17").Select()

$wb.Save()
